$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet: updated mass/weight values ---
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")

$wsGlobal.Range("C6").Value = 24956.015437694587
$wsGlobal.Range("C7").Value = 24560.015437694587
$wsGlobal.Range("C8").Value = 22460.413893925128

$wsGlobal.Range("C12").Value = 20001.523657753718
$wsGlobal.Range("C13").Value = 19605.523657753718
$wsGlobal.Range("C14").Value = 12873.523657753718
$wsGlobal.Range("C15").Value = 11643.980012253716
$wsGlobal.Range("C16").Value = 12058.024012253722

$wsGlobal.Range("C20").Value = 244734.90879206755
$wsGlobal.Range("C21").Value = 240851.47539206757
$wsGlobal.Range("C22").Value = 220261.4179128608

$wsGlobal.Range("C26").Value = 196147.94197831047
$wsGlobal.Range("C27").Value = 192264.50857831046
$wsGlobal.Range("C28").Value = 126246.14077831048
$wsGlobal.Range("C29").Value = 114188.43658716789
$wsGlobal.Range("C30").Value = 118248.82117976793

# --- WING sheet: updated SADRAY / Estimated Mass values ---
$wsWing = $wb.Worksheets.Item("WING")

$wsWing.Range("C10").Value = 2298.0
$wsWing.Range("D10").Value = 14.9

$wsWing.Range("C13").Value = 1892.2857142857142
$wsWing.Range("D13").Value = -5.385714285714288
